$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D21").Value = "[Python] 코사인 유사도(Cosine Similarity)"
$ws.Range("E21").Value = "https://ms-review.tistory.com/18"

$ws.Range("D28").Value = "강화학습 기초"
$ws.Range("E28").Value = "https://ropiens.tistory.com/135"

$ws.Range("D44").Value = "Private 5G의 원리와 동향"
$ws.Range("E44").Value = "https://engineering-ladder.tistory.com/83"
